$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.364.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.80%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.838.71"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.02"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.839.08"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.32%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.483.65"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.845.47"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.429.42"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.43"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.12"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.28%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.36%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.33"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.26"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.989.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.71"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.16%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.804.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +11.72%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.02"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.97"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.317"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.60%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "418.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000295"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.40%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "141.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.55%  "
